# Updates cryptos list data (Price and Volume(1h) columns) per latest scrape
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.988.38'
$ws.Range("E2").Value = '  +1.85%  '
$ws.Range("D3").Value = '1.647.74'
$ws.Range("E3").Value = '  +1.76%  '
$ws.Range("E4").Value = '  -0.14%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '213.62'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.21%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.528'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.63%  '
$ws.Range("E7").Value = '  -0.14%  '
$ws.Range("E9").Value = '  +1.15%  '
$ws.Range("E10").Value = '  +0.23%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0873'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -1.55%  '
$ws.Range("D12").Value = '1.880.23'
$ws.Range("E12").Value = '  +1.69%  '
$ws.Range("D13").Value = '1.645.75'
$ws.Range("E13").Value = '  +1.62%  '
$ws.Range("E14").Value = '  +1.27%  '
$ws.Range("E15").Value = '  +2.81%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '65.63'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +0.73%  '
$ws.Range("D17").Value = '27.974.97'
$ws.Range("E17").Value = '  +1.79%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '233.13'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +1.09%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.69'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +2.15%  '
$ws.Range("D20").Value = '0.0₃0723'
$ws.Range("E20").Value = '  +0.39%  '
$ws.Range("E21").Value = '  -0.09%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '10.69'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +5.20%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '4.40'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +2.68%  '
$ws.Range("E24").Value = '  +4.10%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '152.66'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +1.58%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '6.92'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +1.19%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '15.76'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +1.23%  '
$ws.Range("E28").Value = '  +0.19%  '
$ws.Range("E29").Value = '  -0.13%  '
$ws.Range("E30").Value = '  +1.32%  '
$ws.Range("E31").Value = '  -0.05%  '
$ws.Range("E32").Value = '  +2.78%  '
$ws.Range("D33").Value = '1.449.65'
$ws.Range("E33").Value = '  +0.01%  '
$ws.Range("E34").Value = '  +1.05%  '
$ws.Range("E35").Value = '  +1.75%  '
$ws.Range("E36").Value = '  -0.47%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.890'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +2.97%  '
$ws.Range("E38").Value = '  +0.93%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.562'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +0.19%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.922'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -1.73%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '69.47'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +2.20%  '
$ws.Range("E42").Value = '  +3.02%  '
$ws.Range("E43").Value = '  -0.10%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.48'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.99%  '
$ws.Range("E45").Value = '  +0.89%  '
$ws.Range("E46").Value = '  +5.04%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '5.36'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -1.00%  '
$ws.Range("D48").Value = '1.789.57'
$ws.Range("E48").Value = '  +1.54%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '89.01'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +2.76%  '
$ws.Range("E50").Value = '  -0.72%  '
$ws.Range("E51").Value = '  +0.46%  '
